$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New German instructions text for INSTRUCTIONS_SINGLE (row 4, column B = "de")
$deText = @'
Im Folgenden zeigen wir Ihnen eine Reihe von Personennamen und fragen Sie jeweils, ob es sich bei diesen Personen um literarische Autor*innen handelt, also um Autoren und Autorinnen von literarischer Prosa, Gedichten oder Dramen. Dies trifft nur für einige der gelisteten Namen zu.<br/>
Bitte klicken Sie nur dann „ja“ an, wenn Sie <strong>sich sicher sind</strong>, dass es sich um die Namen literarischer Autor*innen handelt. Wenn Sie nicht wissen, ob es sich um eine/n Autor*in handelt, oder wissen, dass es sich nicht um eine/n Autor*in handelt, klicken Sie „nein/weiß nicht“. Bitte raten Sie nicht. Sie haben für jede Antwort <strong>maximal 10 Sekunden Zeit</strong>. Wenn Sie sich innerhalb dieser Zeit nicht entschieden haben, wird automatisch der nächste Name angezeigt.
'@

# New English instructions text for INSTRUCTIONS_SINGLE (row 4, column C = "en")
$enText = @'
In the following, we show you a number of personal names and ask you in each case whether these persons are literary authors, i.e. authors of literary prose, poetry or drama. This only applies to some of the names listed.<br/> Please only click "yes" if you are <strong>sure</strong> that these are the names of literary authors. If you do not know if it is an author or know that it is not an author, click "no/don't know". Please do not guess. You have <strong>a maximum of 10 seconds</strong> for each answer. If you have not made a decision within this time, the next name will be displayed automatically.
'@

# Replace the two obsolete single-name instructions cells (row 4) with the new, longer
# texts (the two old/duplicate strings that used to live at B4 and C3/C4 are dropped
# because they become unreferenced).
$ws.Range("B4").Value = $deText
$ws.Range("C4").Value = $enText

# Widen column A to fit the key labels, and move the active selection to C4
# (matching the row that now holds the updated instructions text).
$ws.Columns.Item(1).ColumnWidth = 35.0221354166667
$ws.Range("C4").Select()
